# Apply the "Add 2022-06-03 data" update to the Fonds de solidarite workbook.
# Columns: A=dispositif, B=volet, C=nombre_aides, D=nombre_entreprises,
#          E=montant_total, F=reg, G=libelle_region,
#          H=code_categorie_juridique, I=libelle_categorie_juridique
# Only columns C (nombre_aides) and E (montant_total) change for the rows below.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{ Row = 2;   C = 766329;  E = 1429227365 },
    @{ Row = 93;  C = 16941;   E = 50660115 },
    @{ Row = 100; C = 9344;    E = 23829850 },
    @{ Row = 115; C = 81805;   E = 436644943 },
    @{ Row = 121; C = 1306282; E = 2275122346 },
    @{ Row = 129; C = 633627;  E = 3432082404 },
    @{ Row = 130; C = 4247;    E = 141297718 },
    @{ Row = 132; C = 585880;  E = 3468923168 },
    @{ Row = 136; C = 26695;   E = 144329346 },
    @{ Row = 178; C = 515885;  E = 891200284 },
    @{ Row = 237; C = 283320;  E = 1438419995 },
    @{ Row = 240; C = 205917;  E = 1069507520 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 3).Value = $u.C
    $ws.Cells.Item($r, 5).Value = $u.E
}

$wb.Save()
